$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Points de fidélité et monnaie locale"
$ws.Range("B3").Value = "Ceux qui émettent les points de fidélité et ceux qui les recoivent"
$ws.Range("C3").Value = "Ceux qui participent souhaitent de la transparence dans la répartition des points de fidélité"
$ws.Range("D3").Value = "Avoir une vue claire et précise des différents points de fidélité accordés"

$ws.Range("B5").Value = "Tout le monde peut consulter les données mais seuls les utilisateurs enregistrés peuvent bénéficier des points de fidélité et de la monnaie associée"

$ws.Range("A7").Value = "Les acteurs peuvent bénéficier de points de fidélité qui seront transformés en monnaie locale"

$ws.Range("B7").Select()
